$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Junio de 2020 a las 18:57"

# Row 4: Estados Unidos
$ws.Cells.Item(4,1).Value = "Estados Unidos"
$ws.Cells.Item(4,2).Value = 2151113
$ws.Cells.Item(4,3).Value = 8889
$ws.Cells.Item(4,4).Value = 856222
$ws.Cells.Item(4,5).Value = 1177277
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 87
$ws.Cells.Item(4,8).Value = 117614

# Row 5: Brasil
$ws.Cells.Item(5,1).Value = "Brasil"
$ws.Cells.Item(5,2).Value = 852785
$ws.Cells.Item(5,3).Value = 1989
$ws.Cells.Item(5,4).Value = 437512
$ws.Cells.Item(5,5).Value = 372436
$ws.Cells.Item(5,6).Value = 0
$ws.Cells.Item(5,7).Value = 46
$ws.Cells.Item(5,8).Value = 42837

# Row 7: India
$ws.Cells.Item(7,1).Value = "India"
$ws.Cells.Item(7,2).Value = 331724
$ws.Cells.Item(7,3).Value = 10098
$ws.Cells.Item(7,4).Value = 168645
$ws.Cells.Item(7,5).Value = 153594
$ws.Cells.Item(7,6).Value = 0
$ws.Cells.Item(7,7).Value = 286
$ws.Cells.Item(7,8).Value = 9485

# Row 10: Italia
$ws.Cells.Item(10,1).Value = "Italia"
$ws.Cells.Item(10,2).Value = 236989
$ws.Cells.Item(10,3).Value = 338
$ws.Cells.Item(10,4).Value = 176370
$ws.Cells.Item(10,5).Value = 26274
$ws.Cells.Item(10,6).Value = 0
$ws.Cells.Item(10,7).Value = 44
$ws.Cells.Item(10,8).Value = 34345

# Row 12: Alemania
$ws.Cells.Item(12,1).Value = "Alemania"
$ws.Cells.Item(12,2).Value = 187580
$ws.Cells.Item(12,3).Value = 157
$ws.Cells.Item(12,4).Value = 172200
$ws.Cells.Item(12,5).Value = 6513
$ws.Cells.Item(12,6).Value = 0
$ws.Cells.Item(12,7).Value = 0
$ws.Cells.Item(12,8).Value = 8867

# Row 20: Canada
$ws.Cells.Item(20,1).Value = "Canada"
$ws.Cells.Item(20,2).Value = 98735
$ws.Cells.Item(20,3).Value = 325
$ws.Cells.Item(20,4).Value = 60241
$ws.Cells.Item(20,5).Value = 30348
$ws.Cells.Item(20,6).Value = 0
$ws.Cells.Item(20,7).Value = 39
$ws.Cells.Item(20,8).Value = 8146

# Row 42: Irlanda
$ws.Cells.Item(42,1).Value = "Irlanda"
$ws.Cells.Item(42,2).Value = 25303
$ws.Cells.Item(42,3).Value = 8
$ws.Cells.Item(42,4).Value = 22698
$ws.Cells.Item(42,5).Value = 899
$ws.Cells.Item(42,6).Value = 0
$ws.Cells.Item(42,7).Value = 1
$ws.Cells.Item(42,8).Value = 1706

# Row 56: Kazajistan
$ws.Cells.Item(56,1).Value = "Kazajistan"
$ws.Cells.Item(56,2).Value = 14496
$ws.Cells.Item(56,3).Value = 258
$ws.Cells.Item(56,4).Value = 9188
$ws.Cells.Item(56,5).Value = 5231
$ws.Cells.Item(56,6).Value = 0
$ws.Cells.Item(56,7).Value = 4
$ws.Cells.Item(56,8).Value = 77

# Row 62: Argelia
$ws.Cells.Item(62,1).Value = "Argelia"
$ws.Cells.Item(62,2).Value = 10919
$ws.Cells.Item(62,3).Value = 109
$ws.Cells.Item(62,4).Value = 7606
$ws.Cells.Item(62,5).Value = 2546
$ws.Cells.Item(62,6).Value = 0
$ws.Cells.Item(62,7).Value = 7
$ws.Cells.Item(62,8).Value = 767

# Row 76: Uzbekistan
$ws.Cells.Item(76,1).Value = "Uzbekistan"
$ws.Cells.Item(76,2).Value = 5051
$ws.Cells.Item(76,3).Value = 85
$ws.Cells.Item(76,4).Value = 3943
$ws.Cells.Item(76,5).Value = 1089
$ws.Cells.Item(76,6).Value = 0
$ws.Cells.Item(76,7).Value = 0
$ws.Cells.Item(76,8).Value = 19

# Row 87: Kenia
$ws.Cells.Item(87,1).Value = "Kenia"
$ws.Cells.Item(87,2).Value = 3594
$ws.Cells.Item(87,3).Value = 137
$ws.Cells.Item(87,4).Value = 1253
$ws.Cells.Item(87,5).Value = 2238
$ws.Cells.Item(87,6).Value = 0
$ws.Cells.Item(87,7).Value = 3
$ws.Cells.Item(87,8).Value = 103

# Row 95: Somalia
$ws.Cells.Item(95,1).Value = "Somalia"
$ws.Cells.Item(95,2).Value = 2618
$ws.Cells.Item(95,3).Value = 39
$ws.Cells.Item(95,4).Value = 577
$ws.Cells.Item(95,5).Value = 1953
$ws.Cells.Item(95,6).Value = 0
$ws.Cells.Item(95,7).Value = 1
$ws.Cells.Item(95,8).Value = 88

# Row 105: Mali
$ws.Cells.Item(105,1).Value = "Mali"
$ws.Cells.Item(105,2).Value = 1809
$ws.Cells.Item(105,3).Value = 33
$ws.Cells.Item(105,4).Value = 1088
$ws.Cells.Item(105,5).Value = 617
$ws.Cells.Item(105,6).Value = 0
$ws.Cells.Item(105,7).Value = 0
$ws.Cells.Item(105,8).Value = 104

# Row 132: Republica del Chad
$ws.Cells.Item(132,1).Value = "Republica del Chad"
$ws.Cells.Item(132,2).Value = 850
$ws.Cells.Item(132,3).Value = 2
$ws.Cells.Item(132,4).Value = 720
$ws.Cells.Item(132,5).Value = 57
$ws.Cells.Item(132,6).Value = 0
$ws.Cells.Item(132,7).Value = 1
$ws.Cells.Item(132,8).Value = 73

# Row 135: Yemen
$ws.Cells.Item(135,1).Value = "Yemen"
$ws.Cells.Item(135,2).Value = 728
$ws.Cells.Item(135,3).Value = 23
$ws.Cells.Item(135,4).Value = 53
$ws.Cells.Item(135,5).Value = 511
$ws.Cells.Item(135,6).Value = 0
$ws.Cells.Item(135,7).Value = 4
$ws.Cells.Item(135,8).Value = 164

# Row 136: Congo
$ws.Cells.Item(136,1).Value = "Congo"
$ws.Cells.Item(136,2).Value = 728
$ws.Cells.Item(136,3).Value = 0
$ws.Cells.Item(136,4).Value = 221
$ws.Cells.Item(136,5).Value = 483
$ws.Cells.Item(136,6).Value = 0
$ws.Cells.Item(136,7).Value = 0
$ws.Cells.Item(136,8).Value = 24

# Row 137: Crucero
$ws.Cells.Item(137,1).Value = "Crucero"
$ws.Cells.Item(137,2).Value = 712
$ws.Cells.Item(137,3).Value = 0
$ws.Cells.Item(137,4).Value = 651
$ws.Cells.Item(137,5).Value = 48
$ws.Cells.Item(137,6).Value = 0
$ws.Cells.Item(137,7).Value = 13
$ws.Cells.Item(137,8).Value = 13

# Row 151: Liberia
$ws.Cells.Item(151,1).Value = "Liberia"
$ws.Cells.Item(151,2).Value = 458
$ws.Cells.Item(151,3).Value = 12
$ws.Cells.Item(151,4).Value = 219
$ws.Cells.Item(151,5).Value = 207
$ws.Cells.Item(151,6).Value = 0
$ws.Cells.Item(151,7).Value = 0
$ws.Cells.Item(151,8).Value = 32

# Row 206: Groenlandia
$ws.Cells.Item(206,1).Value = "Groenlandia"
$ws.Cells.Item(206,2).Value = 13
$ws.Cells.Item(206,3).Value = 0
$ws.Cells.Item(206,4).Value = 13
$ws.Cells.Item(206,5).Value = 0
$ws.Cells.Item(206,6).Value = 0
$ws.Cells.Item(206,7).Value = 0
$ws.Cells.Item(206,8).Value = 0

# Row 207: Islas Malvinas
$ws.Cells.Item(207,1).Value = "Islas Malvinas"
$ws.Cells.Item(207,2).Value = 13
$ws.Cells.Item(207,3).Value = 0
$ws.Cells.Item(207,4).Value = 13
$ws.Cells.Item(207,5).Value = 0
$ws.Cells.Item(207,6).Value = 0
$ws.Cells.Item(207,7).Value = 0
$ws.Cells.Item(207,8).Value = 0

# Row 208: Islas Turcas y Caicos
$ws.Cells.Item(208,1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(208,2).Value = 12
$ws.Cells.Item(208,3).Value = 0
$ws.Cells.Item(208,4).Value = 11
$ws.Cells.Item(208,5).Value = 0
$ws.Cells.Item(208,6).Value = 0
$ws.Cells.Item(208,7).Value = 0
$ws.Cells.Item(208,8).Value = 1

# Row 209: Santa Sede
$ws.Cells.Item(209,1).Value = "Santa Sede"
$ws.Cells.Item(209,2).Value = 12
$ws.Cells.Item(209,3).Value = 0
$ws.Cells.Item(209,4).Value = 12
$ws.Cells.Item(209,5).Value = 0
$ws.Cells.Item(209,6).Value = 0
$ws.Cells.Item(209,7).Value = 0
$ws.Cells.Item(209,8).Value = 0

# Row 213: Papua Nueva Guinea
$ws.Cells.Item(213,1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(213,2).Value = 8
$ws.Cells.Item(213,3).Value = 0
$ws.Cells.Item(213,4).Value = 8
$ws.Cells.Item(213,5).Value = 0
$ws.Cells.Item(213,6).Value = 0
$ws.Cells.Item(213,7).Value = 0
$ws.Cells.Item(213,8).Value = 0

# Row 214: Islas Virgenes Britanicas
$ws.Cells.Item(214,1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(214,2).Value = 8
$ws.Cells.Item(214,3).Value = 0
$ws.Cells.Item(214,4).Value = 7
$ws.Cells.Item(214,5).Value = 0
$ws.Cells.Item(214,6).Value = 0
$ws.Cells.Item(214,7).Value = 0
$ws.Cells.Item(214,8).Value = 1
